$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.918.74'
$ws.Range("E2").Value = '  +1.07%  '

$ws.Range("D3").Value = '3.963.92'
$ws.Range("E3").Value = '  +3.66%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.733'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.172'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000370'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.80%  '

$ws.Range("D13").Value = '4.597.23'
$ws.Range("E13").Value = '  +4.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").Value = '3.961.59'
$ws.Range("E16").Value = '  +2.57%  '

$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.36%  '

$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").Value = '67.921.53'
$ws.Range("E20").Value = '  +1.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.97%  '

$ws.Range("E23").Value = '  -2.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '732.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.41%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.130'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.28%  '

$ws.Range("E32").Value = '  +2.84%  '

$ws.Range("D33").Value = '0.0₃0892'
$ws.Range("E33").Value = '  +31.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '60.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.30%  '

$ws.Range("E36").Value = '  -4.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.87%  '

$ws.Range("E39").Value = '  -1.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.46%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.33%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.142'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.22%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("E46").Value = '  +1.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.79%  '

$ws.Range("E48").Value = '  -0.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.15%  '
